# Auto-generated edit script: apply scheduled market-data refresh to Seraph_Profits sheets
$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 2278.4285
$ws.Range("I40").Value = 1983.3334
$ws.Range("K40").Value = 1983.3334
$ws.Range("M40").Value = -1808.3334
$ws.Range("H43").Value = 5166
$ws.Range("I43").Value = 1749.5
$ws.Range("J43").Value = 6874.25
$ws.Range("K43").Value = 1749.5
$ws.Range("L43").Value = 6874.25
$ws.Range("M43").Value = -1680.5
$ws.Range("N43").Value = -7012.25
$ws.Range("H62").Value = 7420.4
$ws.Range("I62").Value = 4472.5713
$ws.Range("J62").Value = 9999.75
$ws.Range("K62").Value = 4472.5713
$ws.Range("L62").Value = 9999.75
$ws.Range("M62").Value = -3848.5713
$ws.Range("N62").Value = -11247.75
$ws.Range("H65").Value = 7420.4
$ws.Range("I65").Value = 4472.5713
$ws.Range("J65").Value = 9999.75
$ws.Range("K65").Value = 22362.8565
$ws.Range("L65").Value = 49998.75
$ws.Range("M65").Value = -19242.8565
$ws.Range("N65").Value = -56238.75
$ws.Range("H92").Value = 645.6957
$ws.Range("I92").Value = 549.8823
$ws.Range("J92").Value = 917.1667
$ws.Range("K92").Value = 549.8823
$ws.Range("L92").Value = 917.1667
$ws.Range("M92").Value = 698.1177
$ws.Range("N92").Value = -3413.1667
$ws.Range("H107").Value = 299.125
$ws.Range("I107").Value = 318
$ws.Range("J107").Value = 167
$ws.Range("K107").Value = 318
$ws.Range("L107").Value = 167
$ws.Range("M107").Value = 1602
$ws.Range("N107").Value = -4007
$ws.Range("H138").Value = 4035.4468
$ws.Range("J138").Value = 4275.5713
$ws.Range("L138").Value = 12826.7139
$ws.Range("N138").Value = -23106.7139

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 497.75
$ws.Range("I4").Value = 497.75
$ws.Range("K4").Value = 497.75
$ws.Range("M4").Value = -381.75
$ws.Range("H6").Value = 999999
$ws.Range("I6").Value = 999999
$ws.Range("K6").Value = 999999
$ws.Range("M6").Value = -999826
$ws.Range("H22").Value = 933
$ws.Range("I22").Value = 933
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 933
$ws.Range("L22").Value = 0
$ws.Range("M22").ClearContents()
$ws.Range("N22").Value = -634
$ws.Range("H32").Value = 14916.294
$ws.Range("I32").Value = 14270.19
$ws.Range("J32").Value = 15960
$ws.Range("K32").Value = 14270.19
$ws.Range("L32").Value = 15960
$ws.Range("M32").Value = -13983.19
$ws.Range("N32").Value = -16534
$ws.Range("H74").Value = 1908.4694
$ws.Range("I74").Value = 1220.6136
$ws.Range("J74").Value = 7961.6
$ws.Range("K74").Value = 1220.6136
$ws.Range("L74").Value = 7961.6
$ws.Range("M74").Value = -346.6135999999999
$ws.Range("N74").Value = -9709.6
$ws.Range("H77").Value = 1908.4694
$ws.Range("I77").Value = 1220.6136
$ws.Range("J77").Value = 7961.6
$ws.Range("K77").Value = 6103.067999999999
$ws.Range("L77").Value = 39808
$ws.Range("M77").Value = -1735.067999999999
$ws.Range("N77").Value = -48544
$ws.Range("H122").Value = 2723.1667
$ws.Range("I122").Value = 2126.6316
$ws.Range("K122").Value = 6379.8948
$ws.Range("M122").Value = -3929.8948

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H29").Value = 0
$ws.Range("I29").Value = 0
$ws.Range("K29").Value = 0
$ws.Range("M29").ClearContents()

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 32.6
$ws.Range("I7").Value = 32.6
$ws.Range("K7").Value = 32.6
$ws.Range("M7").Value = 80.40000000000001
$ws.Range("H43").Value = 0
$ws.Range("J43").Value = 0
$ws.Range("L43").ClearContents()
$ws.Range("N43").Value = 0
$ws.Range("H62").Value = 138999.67
$ws.Range("I62").Value = 2000
$ws.Range("K62").Value = 2000
$ws.Range("M62").Value = -1376
$ws.Range("H65").Value = 138999.67
$ws.Range("I65").Value = 2000
$ws.Range("K65").Value = 10000
$ws.Range("M65").Value = -6880
$ws.Range("H99").Value = 15194.357
$ws.Range("I99").Value = 14123.5
$ws.Range("K99").Value = 14123.5
$ws.Range("M99").Value = -12625.5
$ws.Range("H101").Value = 0
$ws.Range("J101").Value = 0
$ws.Range("L101").ClearContents()
$ws.Range("N101").Value = 0
$ws.Range("H103").Value = 9500.75
$ws.Range("I103").Value = 10286.571
$ws.Range("K103").Value = 10286.571
$ws.Range("M103").Value = -9114.571
$ws.Range("H126").Value = 15194.357
$ws.Range("I126").Value = 14123.5
$ws.Range("K126").Value = 42370.5
$ws.Range("M126").Value = -39900.5
$ws.Range("H132").Value = 6998.6665
$ws.Range("I132").Value = 0
$ws.Range("K132").Value = 0
$ws.Range("M132").ClearContents()

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 570.7222
$ws.Range("I5").Value = 619.25
$ws.Range("J5").Value = 556.8570999999999
$ws.Range("K5").Value = 1857.75
$ws.Range("L5").Value = 1670.5713
$ws.Range("M5").Value = -1745.75
$ws.Range("N5").Value = -1894.5713
$ws.Range("H115").Value = 2000
$ws.Range("I115").Value = 2000
$ws.Range("K115").Value = 6000
$ws.Range("M115").Value = -4825
$ws.Range("H128").Value = 3979896.5
$ws.Range("I128").Value = 3979896.5
$ws.Range("K128").Value = 11939689.5
$ws.Range("M128").Value = -11934709.5
$ws.Range("H135").Value = 570.7222
$ws.Range("I135").Value = 619.25
$ws.Range("J135").Value = 556.8570999999999
$ws.Range("K135").Value = 5573.25
$ws.Range("L135").Value = 5011.7139
$ws.Range("M135").Value = -3038.25
$ws.Range("N135").Value = -10081.7139

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H20").Value = 47127
$ws.Range("J20").Value = 61355
$ws.Range("L20").Value = 61355
$ws.Range("N20").Value = -61845
$ws.Range("H24").Value = 22087.75
$ws.Range("J24").Value = 22757.445
$ws.Range("L24").Value = 22757.445
$ws.Range("N24").Value = -23103.445
$ws.Range("H42").Value = 65903.336
$ws.Range("I42").Value = 70420
$ws.Range("J42").Value = 65000
$ws.Range("K42").Value = 70420
$ws.Range("L42").Value = 65000
$ws.Range("M42").Value = -69935
$ws.Range("N42").Value = -65970
$ws.Range("H101").Value = 4913
$ws.Range("J101").Value = 4913
$ws.Range("L101").Value = 4913
$ws.Range("N101").Value = -11403
$ws.Range("H109").Value = 41857.105
$ws.Range("J109").Value = 41857.105
$ws.Range("L109").Value = 41857.105
$ws.Range("N109").Value = -43937.105
$ws.Range("H115").Value = 65903.336
$ws.Range("I115").Value = 70420
$ws.Range("J115").Value = 65000
$ws.Range("K115").Value = 70420
$ws.Range("L115").Value = 65000
$ws.Range("M115").Value = -69245
$ws.Range("N115").Value = -67350
$ws.Range("H122").Value = 57812.555
$ws.Range("I122").Value = 2192.6428
$ws.Range("J122").Value = 252482.25
$ws.Range("K122").Value = 6577.928400000001
$ws.Range("L122").Value = 757446.75
$ws.Range("M122").Value = -4127.928400000001
$ws.Range("N122").Value = -762346.75
$ws.Range("H123").Value = 10000
$ws.Range("J123").Value = 10000
$ws.Range("L123").Value = 10000
$ws.Range("N123").Value = -14900
$ws.Range("H132").Value = 1402.3158
$ws.Range("I132").Value = 1200
$ws.Range("J132").Value = 8888
$ws.Range("K132").Value = 3600
$ws.Range("L132").Value = 26664
$ws.Range("M132").Value = -1070
$ws.Range("N132").Value = -31724

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 0
$ws.Range("I55").Value = 0
$ws.Range("J55").Value = 0
$ws.Range("K55").Value = 0
$ws.Range("L55").ClearContents()
$ws.Range("M55").ClearContents()
$ws.Range("N55").Value = 0
$ws.Range("H97").Value = 15500
$ws.Range("J97").Value = 15500
$ws.Range("L97").Value = 15500
$ws.Range("N97").Value = -17482
$ws.Range("H122").Value = 4875
$ws.Range("I122").Value = 4875
$ws.Range("K122").Value = 14625
$ws.Range("M122").Value = -12175

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H22").Value = 0
$ws.Range("J22").Value = 0
$ws.Range("L22").ClearContents()
$ws.Range("N22").Value = 0
$ws.Range("H44").Value = 60041
$ws.Range("J44").Value = 60041
$ws.Range("L44").Value = 60041
$ws.Range("N44").Value = -61149
$ws.Range("H81").Value = 6105.7334
$ws.Range("I81").Value = 5774.75
$ws.Range("K81").Value = 11549.5
$ws.Range("M81").Value = -10488.5
$ws.Range("H84").Value = 6105.7334
$ws.Range("I84").Value = 5774.75
$ws.Range("K84").Value = 57747.5
$ws.Range("M84").Value = -52443.5
$ws.Range("H107").Value = 416.52
$ws.Range("I107").Value = 291.42856
$ws.Range("J107").Value = 575.7273
$ws.Range("K107").Value = 874.28568
$ws.Range("L107").Value = 1727.1819
$ws.Range("M107").Value = 1045.71432
$ws.Range("N107").Value = -5567.1819
$ws.Range("H113").Value = 267.94116
$ws.Range("I113").Value = 341.1
$ws.Range("J113").Value = 163.42857
$ws.Range("K113").Value = 1023.3
$ws.Range("L113").Value = 490.28571
$ws.Range("M113").Value = 1146.7
$ws.Range("N113").Value = -4830.28571
$ws.Range("H126").Value = 0
$ws.Range("I126").Value = 0
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 0
$ws.Range("L126").ClearContents()
$ws.Range("M126").ClearContents()
$ws.Range("N126").Value = 0
$ws.Range("H132").Value = 1091.25
$ws.Range("I132").Value = 835.7692
$ws.Range("J132").Value = 2198.3333
$ws.Range("K132").Value = 2507.3076
$ws.Range("L132").Value = 6594.999899999999
$ws.Range("M132").Value = 22.69239999999991
$ws.Range("N132").Value = -11654.9999
